# Update "参与人数" (F column) figures across the four sheets of the
# 北京-漫展信息 workbook, as published for the gh-pages data refresh
# (commit 456a3b4). Only numeric values in column F change; everything
# else (G, H, I, styles, text) stays untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 301
$ws.Range("F6").Value  = 433
$ws.Range("F8").Value  = 1960
$ws.Range("F11").Value = 27
$ws.Range("F14").Value = 1314
$ws.Range("F20").Value = 433
$ws.Range("F24").Value = 6949
$ws.Range("F25").Value = 7513
$ws.Range("F26").Value = 32
$ws.Range("F28").Value = 491
$ws.Range("F29").Value = 52
$ws.Range("F31").Value = 238
$ws.Range("F36").Value = 1370
$ws.Range("F37").Value = 10
$ws.Range("F39").Value = 275
$ws.Range("F40").Value = 672
$ws.Range("F43").Value = 307
$ws.Range("F44").Value = 204

# --- Sheet "演出" ---------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 15

# --- Sheet "本地生活" -------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 109

# --- Sheet "全部类型" -------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 15
$ws.Range("F5").Value  = 109
$ws.Range("F7").Value  = 433
$ws.Range("F9").Value  = 1960
$ws.Range("F11").Value = 27
$ws.Range("F15").Value = 1314
$ws.Range("F19").Value = 433
$ws.Range("F23").Value = 6949
$ws.Range("F24").Value = 7513
$ws.Range("F25").Value = 32
$ws.Range("F27").Value = 238
$ws.Range("F30").Value = 10
$ws.Range("F33").Value = 275
$ws.Range("F36").Value = 672
$ws.Range("F42").Value = 307
$ws.Range("F43").Value = 204
